$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Paine, 2021-12-09) was inserted as row 56,
# pushing the existing rows 56-113 down to 57-114.
$ws.Rows.Item(56).Insert()

$ws.Cells.Item(56, 1).Value = 7
$ws.Cells.Item(56, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(56, 3).Value = "Ñuble"
$ws.Cells.Item(56, 4).Value = 44539
$ws.Cells.Item(56, 5).Value = 16
$ws.Cells.Item(56, 6).Value = 100112045
$ws.Cells.Item(56, 7).Value = "Zapallo"
$ws.Cells.Item(56, 8).Value = "Paine"
$ws.Cells.Item(56, 9).Value = "1a (guarda)"
$ws.Cells.Item(56, 10).Value = 120
$ws.Cells.Item(56, 11).Value = 220
$ws.Cells.Item(56, 12).Value = 250
$ws.Cells.Item(56, 13).Value = 235
$ws.Cells.Item(56, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(56, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(56, 16).Value = 235
$ws.Cells.Item(56, 17).Value = 1
$ws.Cells.Item(56, 18).Value = "Hortaliza"
